$d = $word.ActiveDocument

# --- 1. Remove the stray _GoBack bookmark that currently sits in the empty
#        paragraph between the student list and the "Datas e desenvolvimento"
#        heading. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Append a new list item after the last paragraph of the document
#        ("13/08: ... GitHub.") describing the start of the "lista de
#        presidentes". ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = "14/08: Início da lista de presidentes."

# Recompute the paragraph after setting its text (formatting/numbering is
# inherited automatically from the preceding list paragraph).
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$paraRange = $newPara.Range
$textEnd = $paraRange.End - 1

# --- 3. Re-create the _GoBack bookmark collapsed right after the new text
#        (matching Word's usual placement at the last edited position).
#        A directly-collapsed Range confuses Bookmarks.Add in this runtime,
#        so bookmark a one-character placeholder at the insertion point and
#        then delete that character, leaving the bookmark collapsed in
#        place. ---
$placeholderRange = $d.Range($textEnd, $textEnd)
$placeholderRange.InsertAfter("X")

$bookmarkRange = $d.Range($textEnd, $textEnd + 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$deleteRange = $d.Range($textEnd, $textEnd + 1)
$deleteRange.Delete()
